# Update cryptocurrency price/volume figures and fix the Stellar/WhiteBITCoin
# row ordering (rows 50-51 swapped with refreshed data) per the scheduled
# GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.360.82"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "'3.031.07"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'559.15"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'155.59"
$ws.Range("E6").Value = "  -4.14%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.561"
$ws.Range("E8").Value = "  -4.32%  "
$ws.Range("D9").Value = "'3.040.31"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").Value = "'6.42"
$ws.Range("E11").Value = "  -4.53%  "
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "'3.560.57"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").Value = "'63.339.28"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'24.16"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "'3.032.71"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").Value = "'400.08"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "'12.05"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("E22").Value = "  -4.76%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'65.54"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("E25").Value = "  -3.72%  "
$ws.Range("D26").Value = "'0.464"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").Value = "'0.0₃0988"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").Value = "'8.74"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").Value = "'20.45"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").Value = "'162.86"
$ws.Range("E33").Value = "  +7.39%  "
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").Value = "'4.73"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'2.545.83"
$ws.Range("E38").Value = "  -5.61%  "
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "'22.90"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").Value = "'3.96"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("D42").Value = "'37.85"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "'0.671"
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").Value = "'0.0601"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").Value = "'5.12"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "'20.40"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").Value = "'271.39"
$ws.Range("E49").Value = "  -3.79%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.0944"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.48"
$ws.Range("E51").Value = "  +0.41%  "
